# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# on the active worksheet to match the latest scraped values.
#
# Column D ("Price") values that look like plain numbers are written with a
# leading apostrophe (quote-prefix) so Excel keeps them as text, matching the
# original inline-string cell type instead of silently coercing them to
# numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '56.305.19'
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").Value = '2.369.74'
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").Value = '''506.80'
$ws.Range("E5").Value = '  +0.51%  '
$ws.Range("D6").Value = '''130.25'
$ws.Range("E6").Value = '  -1.25%  '
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("E8").Value = '  -1.63%  '
$ws.Range("D9").Value = '2.377.95'
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").Value = '''0.0987'
$ws.Range("E10").Value = '  +1.65%  '
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").Value = '''4.90'
$ws.Range("E12").Value = '  +7.44%  '
$ws.Range("D13").Value = '''0.325'
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").Value = '2.790.07'
$ws.Range("E14").Value = '  -0.94%  '
$ws.Range("D15").Value = '56.288.76'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").Value = '''21.76'
$ws.Range("E16").Value = '  +0.24%  '
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").Value = '2.346.83'
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("D19").Value = '''10.02'
$ws.Range("E19").Value = '  -1.78%  '
$ws.Range("D20").Value = '''310.07'
$ws.Range("E20").Value = '  +0.27%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '''6.27'
$ws.Range("E22").Value = '  +0.26%  '
$ws.Range("D23").Value = '''0.999'
$ws.Range("E23").Value = '  -0.09%  '
$ws.Range("D24").Value = '''65.72'
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  -0.52%  '
$ws.Range("D27").Value = '''0.148'
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("E28").Value = '  -2.95%  '
$ws.Range("D29").Value = '''173.17'
$ws.Range("E29").Value = '  -1.23%  '
$ws.Range("D30").Value = '0.0₃0713'
$ws.Range("E30").Value = '  -1.13%  '
$ws.Range("E31").Value = '  -1.10%  '
$ws.Range("E32").Value = '  -0.14%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").Value = '''0.996'
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -3.49%  '
$ws.Range("E36").Value = '  -1.79%  '
$ws.Range("D37").Value = '''1.19'
$ws.Range("E37").Value = '  -0.52%  '
$ws.Range("D38").Value = '''3.69'
$ws.Range("E38").Value = '  -2.96%  '
$ws.Range("D39").Value = '''0.831'
$ws.Range("E39").Value = '  +1.48%  '
$ws.Range("D40").Value = '''36.36'
$ws.Range("E40").Value = '  -1.33%  '
$ws.Range("E41").Value = '  -3.27%  '
$ws.Range("D42").Value = '''3.38'
$ws.Range("E42").Value = '  +0.65%  '
$ws.Range("D43").Value = '''125.82'
$ws.Range("E43").Value = '  -5.38%  '
$ws.Range("D44").Value = '''4.75'
$ws.Range("D45").Value = '''0.563'
$ws.Range("E45").Value = '  -0.06%  '
$ws.Range("D46").Value = '''0.0899'
$ws.Range("E46").Value = '  -0.98%  '
$ws.Range("D47").Value = '''238.80'
$ws.Range("E47").Value = '  -4.49%  '
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("E49").Value = '  -1.34%  '
$ws.Range("D50").Value = '''16.93'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("E51").Value = '  +0.28%  '
